$wb = $excel.ActiveWorkbook

# --- Sheet: Главные (sheet2 / Worksheets.Item(2)) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("AA2").Value = "2025-11-19 03:11:26"
$ws2.Range("C3").Value = 26
$ws2.Range("D3").Value = 460
$ws2.Range("E3").Value = 207
$ws2.Range("F3").Value = 253
$ws2.Range("G3").Value = 17.69
$ws2.Range("H3").Value = 7.96
$ws2.Range("I3").Value = 9.73
$ws2.Range("J3").Value = 101
$ws2.Range("K3").Value = 104
$ws2.Range("M3").Value = 3
$ws2.Range("AA3").Value = "2025-11-19 03:11:26"
$ws2.Range("C4").Value = 20
$ws2.Range("D4").Value = 326
$ws2.Range("E4").Value = 139
$ws2.Range("F4").Value = 187
$ws2.Range("G4").Value = 16.3
$ws2.Range("H4").Value = 6.95
$ws2.Range("I4").Value = 9.35
$ws2.Range("J4").Value = 67
$ws2.Range("K4").Value = 81
$ws2.Range("AA4").Value = "2025-11-19 03:11:26"
$ws2.Range("C5").Value = 26
$ws2.Range("D5").Value = 423
$ws2.Range("E5").Value = 224
$ws2.Range("F5").Value = 199
$ws2.Range("G5").Value = 16.27
$ws2.Range("H5").Value = 8.619999999999999
$ws2.Range("I5").Value = 7.65
$ws2.Range("J5").Value = 107
$ws2.Range("K5").Value = 92
$ws2.Range("V5").Value = 20
$ws2.Range("AA5").Value = "2025-11-19 03:11:26"
$ws2.Range("C6").Value = 26
$ws2.Range("D6").Value = 447
$ws2.Range("E6").Value = 195
$ws2.Range("F6").Value = 252
$ws2.Range("G6").Value = 17.19
$ws2.Range("H6").Value = 7.5
$ws2.Range("I6").Value = 9.69
$ws2.Range("J6").Value = 90
$ws2.Range("K6").Value = 106
$ws2.Range("AA6").Value = "2025-11-19 03:11:26"
$ws2.Range("C7").Value = 17
$ws2.Range("D7").Value = 225
$ws2.Range("E7").Value = 100
$ws2.Range("F7").Value = 125
$ws2.Range("G7").Value = 13.24
$ws2.Range("I7").Value = 7.35
$ws2.Range("J7").Value = 50
$ws2.Range("K7").Value = 45
$ws2.Range("AA7").Value = "2025-11-19 03:11:26"
$ws2.Range("AA8").Value = "2025-11-19 03:11:26"
$ws2.Range("C9").Value = 26
$ws2.Range("D9").Value = 396
$ws2.Range("E9").Value = 213
$ws2.Range("F9").Value = 183
$ws2.Range("G9").Value = 15.23
$ws2.Range("H9").Value = 8.19
$ws2.Range("I9").Value = 7.04
$ws2.Range("J9").Value = 104
$ws2.Range("K9").Value = 89
$ws2.Range("V9").Value = 14
$ws2.Range("AA9").Value = "2025-11-19 03:11:26"
$ws2.Range("AA10").Value = "2025-11-19 03:11:26"
$ws2.Range("AA11").Value = "2025-11-19 03:11:26"
$ws2.Range("C12").Value = 17
$ws2.Range("D12").Value = 297
$ws2.Range("E12").Value = 131
$ws2.Range("F12").Value = 166
$ws2.Range("G12").Value = 17.47
$ws2.Range("H12").Value = 7.71
$ws2.Range("I12").Value = 9.76
$ws2.Range("J12").Value = 53
$ws2.Range("K12").Value = 58
$ws2.Range("AA12").Value = "2025-11-19 03:11:26"
$ws2.Range("AA13").Value = "2025-11-19 03:11:26"
$ws2.Range("AA14").Value = "2025-11-19 03:11:26"
$ws2.Range("AA15").Value = "2025-11-19 03:11:26"
$ws2.Range("C16").Value = 26
$ws2.Range("D16").Value = 487
$ws2.Range("E16").Value = 242
$ws2.Range("G16").Value = 18.73
$ws2.Range("H16").Value = 9.31
$ws2.Range("I16").Value = 9.42
$ws2.Range("J16").Value = 91
$ws2.Range("X16").Value = 9
$ws2.Range("AA16").Value = "2025-11-19 03:11:26"
$ws2.Range("AA17").Value = "2025-11-19 03:11:26"
$ws2.Range("C18").Value = 25
$ws2.Range("D18").Value = 383
$ws2.Range("E18").Value = 180
$ws2.Range("F18").Value = 203
$ws2.Range("G18").Value = 15.32
$ws2.Range("H18").Value = 7.2
$ws2.Range("I18").Value = 8.119999999999999
$ws2.Range("J18").Value = 80
$ws2.Range("K18").Value = 94
$ws2.Range("M18").Value = 1
$ws2.Range("AA18").Value = "2025-11-19 03:11:26"
$ws2.Range("C19").Value = 21
$ws2.Range("D19").Value = 366
$ws2.Range("E19").Value = 180
$ws2.Range("F19").Value = 186
$ws2.Range("G19").Value = 17.43
$ws2.Range("H19").Value = 8.57
$ws2.Range("I19").Value = 8.859999999999999
$ws2.Range("J19").Value = 85
$ws2.Range("K19").Value = 78
$ws2.Range("AA19").Value = "2025-11-19 03:11:26"
$ws2.Range("AA20").Value = "2025-11-19 03:11:26"
$ws2.Range("C21").Value = 22
$ws2.Range("D21").Value = 312
$ws2.Range("E21").Value = 142
$ws2.Range("G21").Value = 14.18
$ws2.Range("H21").Value = 6.45
$ws2.Range("I21").Value = 7.73
$ws2.Range("J21").Value = 61
$ws2.Range("X21").Value = 2
$ws2.Range("AA21").Value = "2025-11-19 03:11:26"
$ws2.Range("AA22").Value = "2025-11-19 03:11:26"
$ws2.Range("C23").Value = 16
$ws2.Range("D23").Value = 210
$ws2.Range("E23").Value = 77
$ws2.Range("F23").Value = 133
$ws2.Range("G23").Value = 13.13
$ws2.Range("H23").Value = 4.81
$ws2.Range("I23").Value = 8.31
$ws2.Range("J23").Value = 36
$ws2.Range("K23").Value = 54
$ws2.Range("AA23").Value = "2025-11-19 03:11:26"
$ws2.Range("AA24").Value = "2025-11-19 03:11:26"
$ws2.Range("AA25").Value = "2025-11-19 03:11:26"
$ws2.Range("AA26").Value = "2025-11-19 03:11:26"

# --- Sheet: Линейные (sheet3 / Worksheets.Item(3)) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("C2").Value = 17
$ws3.Range("D2").Value = 314
$ws3.Range("E2").Value = 137
$ws3.Range("F2").Value = 177
$ws3.Range("G2").Value = 18.47
$ws3.Range("H2").Value = 8.06
$ws3.Range("I2").Value = 10.41
$ws3.Range("J2").Value = 61
$ws3.Range("K2").Value = 66
$ws3.Range("V2").Value = 8
$ws3.Range("AA2").Value = "2025-11-19 03:11:26"
$ws3.Range("C3").Value = 25
$ws3.Range("D3").Value = 361
$ws3.Range("E3").Value = 183
$ws3.Range("F3").Value = 178
$ws3.Range("G3").Value = 14.44
$ws3.Range("H3").Value = 7.32
$ws3.Range("I3").Value = 7.12
$ws3.Range("J3").Value = 89
$ws3.Range("K3").Value = 74
$ws3.Range("AA3").Value = "2025-11-19 03:11:26"
$ws3.Range("C4").Value = 13
$ws3.Range("D4").Value = 196
$ws3.Range("E4").Value = 88
$ws3.Range("G4").Value = 15.08
$ws3.Range("H4").Value = 6.77
$ws3.Range("I4").Value = 8.31
$ws3.Range("J4").Value = 44
$ws3.Range("X4").Value = 2
$ws3.Range("AA4").Value = "2025-11-19 03:11:26"
$ws3.Range("C5").Value = 13
$ws3.Range("D5").Value = 184
$ws3.Range("E5").Value = 98
$ws3.Range("F5").Value = 86
$ws3.Range("G5").Value = 14.15
$ws3.Range("H5").Value = 7.54
$ws3.Range("I5").Value = 6.62
$ws3.Range("J5").Value = 49
$ws3.Range("K5").Value = 43
$ws3.Range("AA5").Value = "2025-11-19 03:11:26"
$ws3.Range("C6").Value = 16
$ws3.Range("D6").Value = 285
$ws3.Range("E6").Value = 131
$ws3.Range("F6").Value = 154
$ws3.Range("G6").Value = 17.81
$ws3.Range("H6").Value = 8.19
$ws3.Range("I6").Value = 9.630000000000001
$ws3.Range("J6").Value = 58
$ws3.Range("K6").Value = 72
$ws3.Range("AA6").Value = "2025-11-19 03:11:26"
$ws3.Range("C7").Value = 16
$ws3.Range("D7").Value = 253
$ws3.Range("E7").Value = 86
$ws3.Range("F7").Value = 167
$ws3.Range("G7").Value = 15.81
$ws3.Range("H7").Value = 5.38
$ws3.Range("I7").Value = 10.44
$ws3.Range("J7").Value = 43
$ws3.Range("K7").Value = 56
$ws3.Range("AA7").Value = "2025-11-19 03:11:26"
$ws3.Range("C8").Value = 24
$ws3.Range("D8").Value = 379
$ws3.Range("E8").Value = 144
$ws3.Range("F8").Value = 235
$ws3.Range("G8").Value = 15.79
$ws3.Range("H8").Value = 6
$ws3.Range("I8").Value = 9.789999999999999
$ws3.Range("J8").Value = 67
$ws3.Range("K8").Value = 90
$ws3.Range("AA8").Value = "2025-11-19 03:11:26"
$ws3.Range("C9").Value = 24
$ws3.Range("D9").Value = 454
$ws3.Range("E9").Value = 197
$ws3.Range("F9").Value = 257
$ws3.Range("G9").Value = 18.92
$ws3.Range("H9").Value = 8.210000000000001
$ws3.Range("I9").Value = 10.71
$ws3.Range("J9").Value = 86
$ws3.Range("K9").Value = 106
$ws3.Range("AA9").Value = "2025-11-19 03:11:26"
$ws3.Range("AA10").Value = "2025-11-19 03:11:26"
$ws3.Range("AA11").Value = "2025-11-19 03:11:26"
$ws3.Range("AA12").Value = "2025-11-19 03:11:26"
$ws3.Range("AA13").Value = "2025-11-19 03:11:26"
$ws3.Range("AA14").Value = "2025-11-19 03:11:26"
$ws3.Range("C15").Value = 22
$ws3.Range("D15").Value = 425
$ws3.Range("E15").Value = 225
$ws3.Range("F15").Value = 200
$ws3.Range("G15").Value = 19.32
$ws3.Range("H15").Value = 10.23
$ws3.Range("I15").Value = 9.09
$ws3.Range("J15").Value = 90
$ws3.Range("K15").Value = 80
$ws3.Range("V15").Value = 10
$ws3.Range("AA15").Value = "2025-11-19 03:11:26"
$ws3.Range("C16").Value = 25
$ws3.Range("D16").Value = 441
$ws3.Range("E16").Value = 208
$ws3.Range("F16").Value = 233
$ws3.Range("G16").Value = 17.64
$ws3.Range("H16").Value = 8.32
$ws3.Range("I16").Value = 9.32
$ws3.Range("J16").Value = 94
$ws3.Range("K16").Value = 99
$ws3.Range("M16").Value = 7
$ws3.Range("AA16").Value = "2025-11-19 03:11:26"
$ws3.Range("AA17").Value = "2025-11-19 03:11:26"
$ws3.Range("AA18").Value = "2025-11-19 03:11:26"
$ws3.Range("AA19").Value = "2025-11-19 03:11:26"
$ws3.Range("AA20").Value = "2025-11-19 03:11:26"
$ws3.Range("AA21").Value = "2025-11-19 03:11:26"
$ws3.Range("AA22").Value = "2025-11-19 03:11:26"
$ws3.Range("C23").Value = 15
$ws3.Range("D23").Value = 220
$ws3.Range("E23").Value = 107
$ws3.Range("G23").Value = 14.67
$ws3.Range("H23").Value = 7.13
$ws3.Range("I23").Value = 7.53
$ws3.Range("J23").Value = 51
$ws3.Range("X23").Value = 6
$ws3.Range("AA23").Value = "2025-11-19 03:11:26"
$ws3.Range("AA24").Value = "2025-11-19 03:11:26"
$ws3.Range("AA25").Value = "2025-11-19 03:11:26"
$ws3.Range("AA26").Value = "2025-11-19 03:11:26"
